# Scheduled market-data refresh for Sheets workbook (Cactuar server).
# Source data pulled from Universalis; this run refreshes the
# currentAveragePrice(NQ)/LevePrice(NQ/HQ)/LeveProfit(NQ/HQ) columns
# (H, I, J, K, L, M, N) for the affected leve rows on each of the
# eight crafting-class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook
$updated = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8698671
$updated++
$ws.Range("J64").Value = 3773.25
$updated++
$ws.Range("L64").Value = 3773.25
$updated++
$ws.Range("N64").Value = -4269.25
$updated++
$ws.Range("H67").Value = 8698671
$updated++
$ws.Range("J67").Value = 3773.25
$updated++
$ws.Range("L67").Value = 3773.25
$updated++
$ws.Range("N67").Value = -5489.25
$updated++
$ws.Range("H74").Value = 111115830
$updated++
$ws.Range("I74").Value = 500002000
$updated++
$ws.Range("J74").Value = 5499.857
$updated++
$ws.Range("K74").Value = 500002000
$updated++
$ws.Range("L74").Value = 5499.857
$updated++
$ws.Range("M74").Value = -500001064
$updated++
$ws.Range("N74").Value = -7371.857
$updated++
$ws.Range("H77").Value = 111115830
$updated++
$ws.Range("I77").Value = 500002000
$updated++
$ws.Range("J77").Value = 5499.857
$updated++
$ws.Range("K77").Value = 2500010000
$updated++
$ws.Range("L77").Value = 27499.285
$updated++
$ws.Range("M77").Value = -2500005320
$updated++
$ws.Range("N77").Value = -36859.285
$updated++
$ws.Range("H81").Value = 20328
$updated++
$ws.Range("J81").Value = 20328
$updated++
$ws.Range("L81").Value = 20328
$updated++
$ws.Range("N81").Value = -22324
$updated++
$ws.Range("H84").Value = 20328
$updated++
$ws.Range("J84").Value = 20328
$updated++
$ws.Range("L84").Value = 60984
$updated++
$ws.Range("N84").Value = -70968
$updated++
$ws.Range("H86").Value = 1200017.4
$updated++
$ws.Range("I86").Value = 2061079.1
$updated++
$ws.Range("K86").Value = 2061079.1
$updated++
$ws.Range("M86").Value = -2059956.1
$updated++
$ws.Range("H89").Value = 1200017.4
$updated++
$ws.Range("I89").Value = 2061079.1
$updated++
$ws.Range("K89").Value = 10305395.5
$updated++
$ws.Range("M89").Value = -10299779.5
$updated++
$ws.Range("H107").Value = 666
$updated++
$ws.Range("J107").Value = 1000.25
$updated++
$ws.Range("L107").Value = 1000.25
$updated++
$ws.Range("N107").Value = -4840.25
$updated++
$ws.Range("H125").Value = 2834.1177
$updated++
$ws.Range("I125").Value = 1614.8334
$updated++
$ws.Range("K125").Value = 14533.5006
$updated++
$ws.Range("M125").Value = -12073.5006
$updated++
$ws.Range("H137").Value = 9813644
$updated++
$ws.Range("J137").Value = 17549758
$updated++
$ws.Range("L137").Value = 52649274
$updated++
$ws.Range("N137").Value = -52654374
$updated++
$ws.Range("H138").Value = 2993.92
$updated++
$ws.Range("I138").Value = 1347.0435
$updated++
$ws.Range("J138").Value = 3485.8442
$updated++
$ws.Range("K138").Value = 4041.1305
$updated++
$ws.Range("L138").Value = 10457.5326
$updated++
$ws.Range("M138").Value = 1098.8695
$updated++
$ws.Range("N138").Value = -20737.5326
$updated++
$ws.Range("H141").Value = 2921.0293
$updated++
$ws.Range("I141").Value = 2921.0293
$updated++
$ws.Range("J141").Value = 0
$updated++
$ws.Range("K141").Value = 8763.0879
$updated++
$ws.Range("L141").Value = 0
$updated++
$ws.Range("M141").Value = -3583.0879
$updated++
$ws.Range("N141").ClearContents()
$updated++

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 34702.4
$updated++
$ws.Range("I61").Value = 52504
$updated++
$ws.Range("K61").Value = 52504
$updated++
$ws.Range("M61").Value = -52292
$updated++
$ws.Range("H74").Value = 1782.2273
$updated++
$ws.Range("I74").Value = 1401.6154
$updated++
$ws.Range("K74").Value = 1401.6154
$updated++
$ws.Range("M74").Value = -527.6153999999999
$updated++
$ws.Range("H77").Value = 1782.2273
$updated++
$ws.Range("I77").Value = 1401.6154
$updated++
$ws.Range("K77").Value = 7008.076999999999
$updated++
$ws.Range("M77").Value = -2640.076999999999
$updated++
$ws.Range("H136").Value = 34702.4
$updated++
$ws.Range("I136").Value = 52504
$updated++
$ws.Range("K136").Value = 157512
$updated++
$ws.Range("M136").Value = -154962
$updated++

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2861.4878
$updated++
$ws.Range("I20").Value = 2304.0715
$updated++
$ws.Range("K20").Value = 2304.0715
$updated++
$ws.Range("M20").Value = -2057.0715
$updated++
$ws.Range("H105").Value = 2534.4285
$updated++
$ws.Range("I105").Value = 2358.3
$updated++
$ws.Range("K105").Value = 2358.3
$updated++
$ws.Range("M105").Value = -611.3000000000002
$updated++
$ws.Range("H134").Value = 970.8108
$updated++
$ws.Range("I134").Value = 900.55554
$updated++
$ws.Range("J134").Value = 3500
$updated++
$ws.Range("K134").Value = 2701.66662
$updated++
$ws.Range("L134").Value = 10500
$updated++
$ws.Range("M134").Value = -166.66662
$updated++
$ws.Range("N134").Value = -15570
$updated++

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2761.3333
$updated++
$ws.Range("I16").Value = 2695
$updated++
$ws.Range("K16").Value = 2695
$updated++
$ws.Range("M16").Value = -2408
$updated++
$ws.Range("H43").Value = 45512.1
$updated++
$ws.Range("J43").Value = 45512.1
$updated++
$ws.Range("L43").Value = 45512.1
$updated++
$ws.Range("N43").Value = -45880.1
$updated++
$ws.Range("H88").Value = 43791
$updated++
$ws.Range("J88").Value = 43791
$updated++
$ws.Range("L88").Value = 43791
$updated++
$ws.Range("N88").Value = -44603
$updated++
$ws.Range("H91").Value = 43791
$updated++
$ws.Range("J91").Value = 43791
$updated++
$ws.Range("L91").Value = 43791
$updated++
$ws.Range("N91").Value = -46599
$updated++
$ws.Range("H101").Value = 45512.1
$updated++
$ws.Range("J101").Value = 45512.1
$updated++
$ws.Range("L101").Value = 45512.1
$updated++
$ws.Range("N101").Value = -52002.1
$updated++
$ws.Range("H113").Value = 2761.3333
$updated++
$ws.Range("I113").Value = 2695
$updated++
$ws.Range("K113").Value = 2695
$updated++
$ws.Range("M113").Value = -525
$updated++
$ws.Range("H134").Value = 2135.8333
$updated++
$ws.Range("I134").Value = 2143.9412
$updated++
$ws.Range("K134").Value = 6431.823600000001
$updated++
$ws.Range("M134").Value = -3896.823600000001
$updated++

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 215.38095
$updated++
$ws.Range("J38").Value = 278
$updated++
$ws.Range("L38").Value = 834
$updated++
$ws.Range("N38").Value = -1528
$updated++
$ws.Range("H82").Value = 5416.75
$updated++
$ws.Range("I82").Value = 2667
$updated++
$ws.Range("K82").Value = 8001
$updated++
$ws.Range("M82").Value = -7595
$updated++
$ws.Range("H85").Value = 5416.75
$updated++
$ws.Range("I85").Value = 2667
$updated++
$ws.Range("K85").Value = 8001
$updated++
$ws.Range("M85").Value = -6597
$updated++
$ws.Range("H113").Value = 740.55
$updated++
$ws.Range("J113").Value = 774.6667
$updated++
$ws.Range("L113").Value = 2324.0001
$updated++
$ws.Range("N113").Value = -6664.0001
$updated++
$ws.Range("H122").Value = 1117.9231
$updated++
$ws.Range("I122").Value = 769.75
$updated++
$ws.Range("J122").Value = 1272.6666
$updated++
$ws.Range("K122").Value = 6927.75
$updated++
$ws.Range("L122").Value = 11453.9994
$updated++
$ws.Range("M122").Value = -4477.75
$updated++
$ws.Range("N122").Value = -16353.9994
$updated++

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3502241.5
$updated++
$ws.Range("I70").Value = 5055460
$updated++
$ws.Range("J70").Value = 7499.75
$updated++
$ws.Range("K70").Value = 5055460
$updated++
$ws.Range("L70").Value = 7499.75
$updated++
$ws.Range("M70").Value = -5055190
$updated++
$ws.Range("N70").Value = -8039.75
$updated++
$ws.Range("H73").Value = 3502241.5
$updated++
$ws.Range("I73").Value = 5055460
$updated++
$ws.Range("J73").Value = 7499.75
$updated++
$ws.Range("K73").Value = 5055460
$updated++
$ws.Range("L73").Value = 7499.75
$updated++
$ws.Range("M73").Value = -5054524
$updated++
$ws.Range("N73").Value = -9371.75
$updated++
$ws.Range("H97").Value = 676.2381
$updated++
$ws.Range("I97").Value = 526.06665
$updated++
$ws.Range("K97").Value = 526.06665
$updated++
$ws.Range("M97").Value = -30.06664999999998
$updated++
$ws.Range("H102").Value = 20009076
$updated++
$ws.Range("I102").Value = 33342924
$updated++
$ws.Range("J102").Value = 8301.299999999999
$updated++
$ws.Range("K102").Value = 33342924
$updated++
$ws.Range("L102").Value = 8301.299999999999
$updated++
$ws.Range("M102").Value = -33341302
$updated++
$ws.Range("N102").Value = -11545.3
$updated++
$ws.Range("H132").Value = 367715.7
$updated++
$ws.Range("I132").Value = 102496.35
$updated++
$ws.Range("K132").Value = 307489.05
$updated++
$ws.Range("M132").Value = -304959.05
$updated++

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5353.55
$updated++
$ws.Range("I7").Value = 3006.7273
$updated++
$ws.Range("K7").Value = 3006.7273
$updated++
$ws.Range("M7").Value = -2894.7273
$updated++
$ws.Range("H16").Value = 3368.1304
$updated++
$ws.Range("I16").Value = 1830.3125
$updated++
$ws.Range("J16").Value = 6883.143
$updated++
$ws.Range("K16").Value = 1830.3125
$updated++
$ws.Range("L16").Value = 6883.143
$updated++
$ws.Range("M16").Value = -1660.3125
$updated++
$ws.Range("N16").Value = -7223.143
$updated++
$ws.Range("H40").Value = 4099.1665
$updated++
$ws.Range("I40").Value = 4099.1665
$updated++
$ws.Range("K40").Value = 4099.1665
$updated++
$ws.Range("M40").Value = -3963.1665
$updated++
$ws.Range("H104").Value = 33709.855
$updated++
$ws.Range("J104").Value = 33709.855
$updated++
$ws.Range("L104").Value = 33709.855
$updated++
$ws.Range("N104").Value = -40697.855
$updated++
$ws.Range("H122").Value = 12049.363
$updated++
$ws.Range("I122").Value = 4965.8335
$updated++
$ws.Range("K122").Value = 14897.5005
$updated++
$ws.Range("M122").Value = -12447.5005
$updated++
$ws.Range("H126").Value = 5353.55
$updated++
$ws.Range("I126").Value = 3006.7273
$updated++
$ws.Range("K126").Value = 9020.1819
$updated++
$ws.Range("M126").Value = -6550.1819
$updated++

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1751196.6
$updated++
$ws.Range("I81").Value = 3473071
$updated++
$ws.Range("K81").Value = 6946142
$updated++
$ws.Range("M81").Value = -6945081
$updated++
$ws.Range("H84").Value = 1751196.6
$updated++
$ws.Range("I84").Value = 3473071
$updated++
$ws.Range("K84").Value = 34730710
$updated++
$ws.Range("M84").Value = -34725406
$updated++
$ws.Range("H126").Value = 3508.1667
$updated++
$ws.Range("I126").Value = 2809.8
$updated++
$ws.Range("K126").Value = 8429.400000000001
$updated++
$ws.Range("M126").Value = -5959.400000000001
$updated++

Write-Output "Updated $updated cells across 8 sheets"
